$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fmtNum = "#,##0"
$fmtPct = "#,##0.0;`"-`"#,##0.0"

# --- Header text updates (Volume number, week-covering date range) ---
$ws.Range("A8").Value = "Volume 30   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/30/2023  Through  2/5/2023"

# --- Cells changing from a placeholder shared-string ("0") to a real number ---
$ws.Range("C15").NumberFormat = $fmtNum
$ws.Range("C15").Value = 1
$ws.Range("F15").NumberFormat = $fmtNum
$ws.Range("F15").Value = 1
$ws.Range("I15").NumberFormat = $fmtNum
$ws.Range("I15").Value = 1
$ws.Range("C20").NumberFormat = $fmtNum
$ws.Range("C20").Value = 2
$ws.Range("C26").NumberFormat = $fmtNum
$ws.Range("C26").Value = 1
$ws.Range("D28").NumberFormat = $fmtNum
$ws.Range("D28").Value = 3
$ws.Range("E28").NumberFormat = $fmtPct
$ws.Range("E28").Value = -100
$ws.Range("G28").NumberFormat = $fmtNum
$ws.Range("G28").Value = 3
$ws.Range("H28").NumberFormat = $fmtPct
$ws.Range("H28").Value = -33.333333333333
$ws.Range("J28").NumberFormat = $fmtNum
$ws.Range("J28").Value = 3
$ws.Range("K28").NumberFormat = $fmtPct
$ws.Range("K28").Value = 0
$ws.Range("D29").NumberFormat = $fmtNum
$ws.Range("D29").Value = 2
$ws.Range("E29").NumberFormat = $fmtPct
$ws.Range("E29").Value = -100
$ws.Range("G29").NumberFormat = $fmtNum
$ws.Range("G29").Value = 1
$ws.Range("H29").NumberFormat = $fmtPct
$ws.Range("H29").Value = -50
$ws.Range("J29").NumberFormat = $fmtNum
$ws.Range("J29").Value = 2
$ws.Range("K29").NumberFormat = $fmtPct
$ws.Range("K29").Value = 0

# --- Cells changing from a real number back to a placeholder shared-string ---
$ws.Range("D23").Value = "'0"
$ws.Range("E23").Value = "***.*"
$ws.Range("F23").Value = "'0"
$ws.Range("D27").Value = "'0"
$ws.Range("E27").Value = "***.*"

# Re-apply the plain "General" text style (matching other placeholder cells)
# by pasting formats from an existing placeholder cell (A15 uses that exact style).
$ws.Range("A15").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Plain numeric value updates (type unchanged) ---
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -83.333333333333
$ws.Range("C16").Value = 4
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 57.142857142857
$ws.Range("I16").Value = 14
$ws.Range("J16").Value = 12
$ws.Range("K16").Value = 16.666666666666
$ws.Range("L16").Value = 366.666666666667
$ws.Range("M16").Value = -51.724137931034
$ws.Range("N16").Value = -74.074074074074
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -88.888888888888
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 17
$ws.Range("J17").Value = 20
$ws.Range("K17").Value = -15
$ws.Range("L17").Value = -10.526315789473
$ws.Range("M17").Value = 54.545454545454
$ws.Range("N17").Value = -41.379310344827
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 42.857142857142
$ws.Range("I18").Value = 11
$ws.Range("J18").Value = 8
$ws.Range("K18").Value = 37.5
$ws.Range("L18").Value = -35.294117647058
$ws.Range("M18").Value = -57.692307692307
$ws.Range("N18").Value = -89.108910891089
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 900
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 25.806451612903
$ws.Range("I19").Value = 53
$ws.Range("J19").Value = 48
$ws.Range("K19").Value = 10.416666666666
$ws.Range("L19").Value = 96.296296296296
$ws.Range("M19").Value = 12.765957446808
$ws.Range("N19").Value = -10.169491525423
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -55.555555555555
$ws.Range("I20").Value = 7
$ws.Range("J20").Value = 11
$ws.Range("K20").Value = -36.363636363636
$ws.Range("L20").Value = 16.666666666666
$ws.Range("M20").Value = -63.157894736842
$ws.Range("N20").Value = -97.910447761194
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 33.333333333333
$ws.Range("F21").Value = 79
$ws.Range("G21").Value = 72
$ws.Range("H21").Value = 9.722222222222
$ws.Range("I21").Value = 105
$ws.Range("J21").Value = 99
$ws.Range("K21").Value = 6.060606060606
$ws.Range("L21").Value = 45.833333333333
$ws.Range("M21").Value = -21.052631578947
$ws.Range("N21").Value = -82.081911262798
$ws.Range("H23").Value = -100
$ws.Range("M23").Value = -85.714285714285
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 11.764705882352
$ws.Range("F24").Value = 76
$ws.Range("G24").Value = 72
$ws.Range("H24").Value = 5.555555555555
$ws.Range("I24").Value = 103
$ws.Range("J24").Value = 93
$ws.Range("K24").Value = 10.752688172043
$ws.Range("L24").Value = 21.176470588235
$ws.Range("M24").Value = 5.102040816326
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -37.5
$ws.Range("F25").Value = 23
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = 4.545454545454
$ws.Range("I25").Value = 30
$ws.Range("J25").Value = 24
$ws.Range("K25").Value = 25
$ws.Range("L25").Value = 66.666666666666
$ws.Range("M25").Value = -3.225806451612
$ws.Range("F26").Value = 3
$ws.Range("I26").Value = 3
$ws.Range("L26").Value = 200
$ws.Range("I27").Value = 4
$ws.Range("K27").Value = 300
$ws.Range("L27").Value = 300
$ws.Range("M28").Value = -25
$ws.Range("N28").Value = -25
$ws.Range("M29").Value = -33.333333333333
$ws.Range("N29").Value = -50
